$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

$ws.Range("B62").Value = "2025-04-28 05:40:23"

$ws.Range("C62").Value = "James Davis took a picture of New Battery from Ford.`nNow James Davis is Excited, feeling that the task was Stressful.`n"
